$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column B/C are plain text (coin name / link); D/E are
# numeric-looking text (price / volume%) that must stay text, so we
# force a text format before writing and restore the default style after.
$textCells = @{
  "D2" = "54.318.36"
  "E2" = "  -7.64%  "
  "D3" = "2.887.75"
  "E3" = "  -10.64%  "
  "D4" = "0.999"
  "E4" = "  -0.08%  "
  "D5" = "469.13"
  "E5" = "  -12.85%  "
  "D6" = "125.67"
  "E6" = "  -7.88%  "
  "D8" = "2.876.23"
  "E8" = "  -10.96%  "
  "E9" = "  -13.19%  "
  "D10" = "6.52"
  "E10" = "  -14.24%  "
  "D11" = "0.0959"
  "E11" = "  -16.37%  "
  "D12" = "0.328"
  "E12" = "  -16.91%  "
  "E13" = "  -3.94%  "
  "D14" = "3.365.11"
  "E14" = "  -11.21%  "
  "D15" = "22.65"
  "E15" = "  -12.95%  "
  "D16" = "54.361.36"
  "E16" = "  -7.69%  "
  "D17" = "2.877.19"
  "E17" = "  -11.18%  "
  "D18" = "0.0000132"
  "E18" = "  -16.47%  "
  "D19" = "5.07"
  "E19" = "  -14.01%  "
  "D20" = "11.42"
  "E20" = "  -13.94%  "
  "D21" = "7.07"
  "E21" = "  -14.74%  "
  "D22" = "305.05"
  "E22" = "  -15.49%  "
  "E23" = "  +0.09%  "
  "D24" = "0.444"
  "E24" = "  -14.71%  "
  "D25" = "58.73"
  "E25" = "  -16.58%  "
  "D26" = "1.00"
  "E26" = "  +0.35%  "
  "E27" = "  -11.02%  "
  "D28" = "0.997"
  "E28" = "  -0.20%  "
  "D29" = "0.0₃0818"
  "E29" = "  -15.65%  "
  "D30" = "6.02"
  "E30" = "  -14.79%  "
  "D31" = "1.13"
  "E31" = "  -7.96%  "
  "D32" = "6.04"
  "E32" = "  -14.54%  "
  "D33" = "19.02"
  "E33" = "  -13.49%  "
  "D34" = "1.57"
  "E34" = "  -18.77%  "
  "D35" = "140.00"
  "E35" = "  -14.32%  "
  "D36" = "4.18"
  "E36" = "  -15.32%  "
  "D37" = "5.41"
  "E37" = "  -15.52%  "
  "D38" = "1.21"
  "E38" = "  -15.81%  "
  "D39" = "2.895.77"
  "E39" = "  -11.35%  "
  "D40" = "0.0612"
  "E40" = "  -13.79%  "
  "D41" = "0.997"
  "E41" = "  -0.29%  "
  "D42" = "21.66"
  "E42" = "  -17.90%  "
  "D43" = "34.31"
  "E43" = "  -16.55%  "
  "D44" = "0.956"
  "E44" = "  -12.83%  "
  "D45" = "3.41"
  "E45" = "  -15.41%  "
  "D46" = "0.591"
  "E46" = "  -17.57%  "
  "D47" = "1.30"
  "E47" = "  -13.74%  "
  "D48" = "2.035.20"
  "E48" = "  -11.49%  "
  "D49" = "17.85"
  "E49" = "  -14.26%  "
  "D50" = "5.26"
  "E50" = "  -16.53%  "
  "D51" = "0.0208"
  "E51" = "  -14.28%  "
}
$plainCells = @{
  "B32" = "InternetComputer(DFINITY)"
  "C32" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
  "B33" = "EthereumClassic"
  "C33" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
  "B40" = "Hedera"
  "C40" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
  "B41" = "FirstDigitalUSD"
  "C41" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
  "B45" = "Filecoin"
  "C45" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
  "B46" = "Mantle"
  "C46" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
  "B47" = "Stacks"
  "C47" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
  "B48" = "Maker"
  "C48" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
  "B49" = "InjectiveProtocol"
  "C49" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
  "B50" = "Cosmos"
  "C50" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
}

foreach ($ref in $textCells.Keys) {
  $cell = $ws.Range($ref)
  $cell.NumberFormat = "@"
  $cell.Value = $textCells[$ref]
  $cell.Style = "Normal"
}

foreach ($ref in $plainCells.Keys) {
  $ws.Range($ref).Value = $plainCells[$ref]
}
